$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel re-interpreting
# numeric-looking strings (e.g. "1.000", "240.60") as numbers. We briefly
# switch the cell to Text format, assign the value, then restore the
# cell's original number format so no visible formatting changes remain.
function Set-TextValue($rangeAddr, $text) {
    $c = $ws.Range($rangeAddr)
    $origFormat = $c.NumberFormat
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = $origFormat
}

Set-TextValue 'D2' '30.395.17'
Set-TextValue 'E2' '  -0.36%  '
Set-TextValue 'D3' '1.926.31'
Set-TextValue 'E3' '  +3.98%  '
Set-TextValue 'D4' '0.9995'
Set-TextValue 'E4' '  -0.13%  '
Set-TextValue 'D5' '240.60'
Set-TextValue 'E6' '  -0.08%  '
Set-TextValue 'D7' '0.4765'
Set-TextValue 'E7' '  +0.44%  '
Set-TextValue 'D8' '0.2862'
Set-TextValue 'E8' '  +4.27%  '
Set-TextValue 'D9' '0.06589'
Set-TextValue 'E9' '  +4.19%  '
Set-TextValue 'D10' '19.10'
Set-TextValue 'E10' '  +8.25%  '
Set-TextValue 'D11' '106.28'
Set-TextValue 'E11' '  +25.60%  '
Set-TextValue 'D12' '1.917.87'
Set-TextValue 'E12' '  +3.59%  '
Set-TextValue 'D13' '0.07613'
Set-TextValue 'E13' '  +2.24%  '
Set-TextValue 'D14' '5.122'
Set-TextValue 'E14' '  +3.11%  '
Set-TextValue 'D15' '0.6577'
Set-TextValue 'E15' '  +5.22%  '
Set-TextValue 'D16' '302.07'
Set-TextValue 'E16' '  +21.99%  '
Set-TextValue 'D17' '30.401.76'
Set-TextValue 'E17' '  -0.28%  '
Set-TextValue 'E18' '  +0.04%  '
Set-TextValue 'D19' '12.94'
Set-TextValue 'E19' '  +1.79%  '
Set-TextValue 'B20' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C20' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D20' '2.166.85'
Set-TextValue 'E20' '  +2.82%  '
Set-TextValue 'B21' 'ShibaInu'
Set-TextValue 'C21' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D21' '0.000007494'
Set-TextValue 'E21' '  +2.21%  '
Set-TextValue 'E22' '  -0.12%  '
Set-TextValue 'D23' '5.263'
Set-TextValue 'E23' '  +6.85%  '
Set-TextValue 'D24' '6.311'
Set-TextValue 'E24' '  +6.83%  '
Set-TextValue 'D25' '168.60'
Set-TextValue 'E25' '  +2.51%  '
Set-TextValue 'D26' '9.208'
Set-TextValue 'E26' '  +0.93%  '
Set-TextValue 'D27' '19.78'
Set-TextValue 'E27' '  +10.10%  '
Set-TextValue 'D28' '2.001'
Set-TextValue 'E28' '  +7.03%  '
Set-TextValue 'D29' '0.1120'
Set-TextValue 'E29' '  +9.34%  '
Set-TextValue 'E30' '  -0.52%  '
Set-TextValue 'D31' '4.088'
Set-TextValue 'E31' '  +1.16%  '
Set-TextValue 'D32' '3.922'
Set-TextValue 'E32' '  +2.38%  '
Set-TextValue 'D33' '0.05002'
Set-TextValue 'E33' '  +3.21%  '
Set-TextValue 'D34' '0.7403'
Set-TextValue 'E34' '  +6.22%  '
Set-TextValue 'D35' '1.147'
Set-TextValue 'E35' '  +1.34%  '
Set-TextValue 'D36' '0.9996'
Set-TextValue 'E36' '  -0.02%  '
Set-TextValue 'E37' '  +1.07%  '
Set-TextValue 'D38' '0.01947'
Set-TextValue 'E38' '  +2.64%  '
Set-TextValue 'D39' '2.700'
Set-TextValue 'E39' '  +0.61%  '
Set-TextValue 'D40' '2.054'
Set-TextValue 'E40' '  +2.22%  '
Set-TextValue 'D41' '0.8727'
Set-TextValue 'E41' '  -0.21%  '
Set-TextValue 'D42' '107.33'
Set-TextValue 'E42' '  +0.91%  '
Set-TextValue 'D43' '5.782'
Set-TextValue 'E43' '  +4.42%  '
Set-TextValue 'B44' 'PaxDollar'
Set-TextValue 'C44' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D44' '1.000'
Set-TextValue 'E44' '  -0.08%  '
Set-TextValue 'B45' 'Aave'
Set-TextValue 'C45' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D45' '69.63'
Set-TextValue 'E45' '  +10.43%  '
Set-TextValue 'D46' '0.4123'
Set-TextValue 'E46' '  +1.51%  '
Set-TextValue 'D47' '7.210'
Set-TextValue 'E47' '  +0.36%  '
Set-TextValue 'D48' '9.257'
Set-TextValue 'E48' '  +8.34%  '
Set-TextValue 'E49' '  +3.41%  '
Set-TextValue 'E50' '  +0.33%  '
Set-TextValue 'D51' '0.05623'
Set-TextValue 'E51' '  +1.87%  '
